$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ---- Row 14: 2889. Reshape Data: Pivot ----
$row14 = $tbl.ListRows.Add()
$ws.Range("A14").Value = "2889. Reshape Data: Pivot"
$ws.Hyperlinks.Add($ws.Range("E14"), "https://leetcode.com/problems/reshape-data-pivot/solutions/4141174/pandas-1-line-elegant-short-and-more-pandas-solutions/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata") | Out-Null
$ws.Range("E14").Style = $ws.Range("E13").Style
$ws.Range("D14").Value = "Use df.pivot(index='month', columns='city', values='temperature')"
$ws.Range("B14").Value = "Easy"
$ws.Range("B14").Interior.Color = $ws.Range("B13").Interior.Color
$ws.Range("C14").Value = "Table Reshaping"

# ---- Row 15: 2890. Reshape Data: Melt ----
$row15 = $tbl.ListRows.Add()
$ws.Range("A15").Value = "2890. Reshape Data: Melt"
$ws.Range("D15").Value = "Use pd.melt(report, id_vars=['product'], var_name='quarter', value_name='sales')"
$ws.Hyperlinks.Add($ws.Range("E15"), "https://leetcode.com/problems/reshape-data-melt/solutions/4141084/line-by-line-explanation-easy-solution-beginner-friendly-pandas/?envType=study-plan-v2&envId=introduction-to-pandas&lang=pythondata ") | Out-Null
$ws.Range("E15").Style = $ws.Range("E13").Style
$ws.Range("B15").Value = "Easy"
$ws.Range("B15").Interior.Color = $ws.Range("B13").Interior.Color
$ws.Range("C15").Value = "Table Reshaping"

# ---- Row 16: 2891. Method Chaining ----
$row16 = $tbl.ListRows.Add()
$ws.Range("C16").Value = "Advanced Techniques"
$ws.Range("A16").Value = "2891. Method Chaining"
$ws.Range("D16").Value = "Chain methods: return animals[animals['weight'] > 100].sort_values(['weight'], ascending=False,)[['name']]"
$ws.Hyperlinks.Add($ws.Range("E16"), "https://leetcode.com/problems/method-chaining/solutions/4134716/easy-pandas-solution-one-liner-beginner-friendly/ ") | Out-Null
$ws.Range("E16").Style = $ws.Range("E13").Style
$ws.Range("B16").Value = "Easy"
$ws.Range("B16").Interior.Color = $ws.Range("B13").Interior.Color
